$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 3 ----
$ws.Range("A3").Value = 70752145
$ws.Range("B3").Value = 96355
$ws.Range("C3").Value = "Ovaliderad"
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 219862
$ws.Range("F3").Value = "Nästrot"
$ws.Range("G3").Value = "Neottia nidus-avis"
$ws.Range("H3").Value = "(L.) Rich."
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = ""
$ws.Range("L3").Value = ""
$ws.Range("P3").Value = "Myrorna, N om, Upl"
$ws.Range("Q3").Value = 703256.005931965
$ws.Range("R3").Value = 6626085.902598036
$ws.Range("S3").Value = 5
$ws.Range("T3").Value = "Stockholm"
$ws.Range("U3").Value = "Norrtälje"
$ws.Range("V3").Value = "Uppland"
$ws.Range("W3").Value = "Frötuna"
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2018-04-18"
$ws.Range("Y3").Style = "Normal"
$ws.Range("Z3").Value = "00:00"
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "2018-04-18"
$ws.Range("AA3").Style = "Normal"
$ws.Range("AB3").Value = "00:00"
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false
$ws.Range("AT3").Value = ""
$ws.Range("AW3").Value = "Bo Törnquist"
$ws.Range("AX3").Value = "Bo Törnquist, Kjell  Andersson"
$ws.Range("AY3").Value = ""

# ---- Row 4 ----
$ws.Range("AC4").Value = ""
$ws.Range("AF4").Value = ""
$ws.Range("A4").Value = 70744400
$ws.Range("B4").Value = 96334
$ws.Range("C4").Value = "Ovaliderad"
$ws.Range("D4").Value = "VU"
$ws.Range("E4").Value = 220787
$ws.Range("F4").Value = "Knärot"
$ws.Range("G4").Value = "Goodyera repens"
$ws.Range("H4").Value = "(L.) R. Br."
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = ""
$ws.Range("K4").Value = ""
$ws.Range("L4").Value = ""
$ws.Range("P4").Value = "Myrorna, N om, Upl"
$ws.Range("Q4").Value = 703206.0057106519
$ws.Range("R4").Value = 6626260.055830983
$ws.Range("S4").Value = 5
$ws.Range("T4").Value = "Stockholm"
$ws.Range("U4").Value = "Norrtälje"
$ws.Range("V4").Value = "Uppland"
$ws.Range("W4").Value = "Frötuna"
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = "2018-04-18"
$ws.Range("Y4").Style = "Normal"
$ws.Range("Z4").Value = "00:00"
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = "2018-04-18"
$ws.Range("AA4").Style = "Normal"
$ws.Range("AB4").Value = "00:00"
$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AG4").Value = $false
$ws.Range("AT4").Value = ""
$ws.Range("AW4").Value = "Bo Törnquist"
$ws.Range("AX4").Value = "Bo Törnquist, Kjell  Andersson"
$ws.Range("AY4").Value = ""

# ---- Row 5 ----
$ws.Range("L5").Value = ""
$ws.Range("A5").Value = 70743844
$ws.Range("B5").Value = 89410
$ws.Range("C5").Value = "Ovaliderad"
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 5432
$ws.Range("F5").Value = "Granticka"
$ws.Range("G5").Value = "Porodaedalea chrysoloma"
$ws.Range("H5").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = ""
$ws.Range("P5").Value = "Myrorna, N om, Upl"
$ws.Range("Q5").Value = 703367.1782205966
$ws.Range("R5").Value = 6626265.819183424
$ws.Range("S5").Value = 5
$ws.Range("T5").Value = "Stockholm"
$ws.Range("U5").Value = "Norrtälje"
$ws.Range("V5").Value = "Uppland"
$ws.Range("W5").Value = "Frötuna"
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value = "2018-04-18"
$ws.Range("Y5").Style = "Normal"
$ws.Range("Z5").Value = "00:00"
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value = "2018-04-18"
$ws.Range("AA5").Style = "Normal"
$ws.Range("AB5").Value = "00:00"
$ws.Range("AC5").Value = "På grenar på levande gammal gran."
$ws.Range("AD5").Value = $false
$ws.Range("AE5").Value = $false
$ws.Range("AF5").Value = ""
$ws.Range("AG5").Value = $false
$ws.Range("AT5").Value = ""
$ws.Range("AW5").Value = "Bo Törnquist"
$ws.Range("AX5").Value = "Bo Törnquist, Kjell  Andersson"
$ws.Range("AY5").Value = ""

# ---- Row 6 ----
$ws.Range("X6").Value = ""
$ws.Range("A6").Value = 70752141
$ws.Range("B6").Value = 89410
$ws.Range("C6").Value = "Ovaliderad"
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 5432
$ws.Range("F6").Value = "Granticka"
$ws.Range("G6").Value = "Porodaedalea chrysoloma"
$ws.Range("H6").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("I6").Value = ""
$ws.Range("J6").Value = ""
$ws.Range("K6").Value = ""
$ws.Range("P6").Value = "Myrorna, N om, Upl"
$ws.Range("Q6").Value = 703256.005931965
$ws.Range("R6").Value = 6626085.902598036
$ws.Range("S6").Value = 5
$ws.Range("T6").Value = "Stockholm"
$ws.Range("U6").Value = "Norrtälje"
$ws.Range("V6").Value = "Uppland"
$ws.Range("W6").Value = "Frötuna"
$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value = "2018-04-18"
$ws.Range("Y6").Style = "Normal"
$ws.Range("Z6").Value = "00:00"
$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value = "2018-04-18"
$ws.Range("AA6").Style = "Normal"
$ws.Range("AB6").Value = "00:00"
$ws.Range("AC6").Value = "På högstubbe av gran. Gammal barrblandskog. Inslag av försumpade fuktdråg och kalkpåverkad mark."
$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AF6").Value = ""
$ws.Range("AG6").Value = $false
$ws.Range("AT6").Value = ""
$ws.Range("AW6").Value = "Bo Törnquist"
$ws.Range("AX6").Value = "Bo Törnquist, Kjell  Andersson"
$ws.Range("AY6").Value = ""

# ---- Row 7 ----
$ws.Range("A7").Value = 70744212
$ws.Range("B7").Value = 103265
$ws.Range("C7").Value = "Ovaliderad"
$ws.Range("D7").Value = "LC"
$ws.Range("E7").Value = 221144
$ws.Range("F7").Value = "Grönpyrola"
$ws.Range("G7").Value = "Pyrola chlorantha"
$ws.Range("H7").Value = "Sw."
$ws.Range("I7").Value = ""
$ws.Range("J7").Value = ""
$ws.Range("K7").Value = ""
$ws.Range("L7").Value = ""
$ws.Range("P7").Value = "Myrorna, N om, Upl"
$ws.Range("Q7").Value = 703273.8314362012
$ws.Range("R7").Value = 6626260.218834675
$ws.Range("S7").Value = 5
$ws.Range("T7").Value = "Stockholm"
$ws.Range("U7").Value = "Norrtälje"
$ws.Range("V7").Value = "Uppland"
$ws.Range("W7").Value = "Frötuna"
$ws.Range("Y7").NumberFormat = "@"
$ws.Range("Y7").Value = "2018-04-18"
$ws.Range("Y7").Style = "Normal"
$ws.Range("Z7").Value = "00:00"
$ws.Range("AA7").NumberFormat = "@"
$ws.Range("AA7").Value = "2018-04-18"
$ws.Range("AA7").Style = "Normal"
$ws.Range("AB7").Value = "00:00"
$ws.Range("AD7").Value = $false
$ws.Range("AE7").Value = $false
$ws.Range("AG7").Value = $false
$ws.Range("AT7").Value = ""
$ws.Range("AW7").Value = "Bo Törnquist"
$ws.Range("AX7").Value = "Bo Törnquist, Kjell  Andersson"
$ws.Range("AY7").Value = ""

# ---- Row 8 ----
$ws.Range("AF8").Value = ""
$ws.Range("J8").Value = ""
$ws.Range("K8").Value = ""
$ws.Range("A8").Value = 86555940
$ws.Range("B8").Value = 96334
$ws.Range("C8").Value = "Ovaliderad"
$ws.Range("D8").Value = "VU"
$ws.Range("E8").Value = 220787
$ws.Range("F8").Value = "Knärot"
$ws.Range("G8").Value = "Goodyera repens"
$ws.Range("H8").Value = "(L.) R. Br."
$ws.Range("I8").Value = ""
$ws.Range("P8").Value = "Myrorna, N om (*knärot* /stjälk/), Upl"
$ws.Range("Q8").Value = 703206.7403017445
$ws.Range("R8").Value = 6626265.149522305
$ws.Range("S8").Value = 10
$ws.Range("T8").Value = "Stockholm"
$ws.Range("U8").Value = "Norrtälje"
$ws.Range("V8").Value = "Uppland"
$ws.Range("W8").Value = "Frötuna"
$ws.Range("X8").Value = "AB-Nor-1742"
$ws.Range("Y8").NumberFormat = "@"
$ws.Range("Y8").Value = "2018-04-18"
$ws.Range("Y8").Style = "Normal"
$ws.Range("Z8").Value = "00:00"
$ws.Range("AA8").NumberFormat = "@"
$ws.Range("AA8").Value = "2018-04-18"
$ws.Range("AA8").Style = "Normal"
$ws.Range("AB8").Value = "00:00"
$ws.Range("AC8").Value = "Obs: Flera dellokaler. Se privata obsar! Myrorna, N om, Obskoord: 6626090/1658024/10 m (). Enstaka fynd. Gammal barrblandskog. Delvis fuktdråg och kalkpåverkad mark."
$ws.Range("AD8").Value = $false
$ws.Range("AE8").Value = $false
$ws.Range("AG8").Value = $false
$ws.Range("AT8").Value = ""
$ws.Range("AW8").Value = "Jan Yngve Andersson"
$ws.Range("AX8").Value = "Kjell  Andersson, Bo Törnquist"
$ws.Range("AY8").Value = "Floraväkteri Sverige"
